$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = '27.681.05'
$ws.Range("E2").Value = '  +0.62%  '
$ws.Range("D3").Value = '1.850.05'
$ws.Range("E3").Value = '  +0.61%  '
$ws.Range("E4").Value = '  +0.16%  '
$ws.Range("D5").Value = '312.96'
$ws.Range("E5").Value = '  +0.18%  '
$ws.Range("E6").Value = '  +0.16%  '
$ws.Range("D7").Value = '0.4266'
$ws.Range("E7").Value = '  +0.74%  '
$ws.Range("D8").Value = '0.3633'
$ws.Range("E8").Value = '  +0.38%  '
$ws.Range("D9").Value = '44.72'
$ws.Range("E9").Value = '  +2.71%  '
$ws.Range("D10").Value = '0.07300'
$ws.Range("E10").Value = '  +1.63%  '
$ws.Range("D11").Value = '0.8752'
$ws.Range("E11").Value = '  -2.26%  '
$ws.Range("D12").Value = '20.60'
$ws.Range("E12").Value = '  +0.07%  '
$ws.Range("D13").Value = '1.894.45'
$ws.Range("E13").Value = '  +3.81%  '
$ws.Range("D14").Value = '5.323'
$ws.Range("E14").Value = '  +0.30%  '
$ws.Range("D15").Value = '6.510'
$ws.Range("E15").Value = '  -1.08%  '
$ws.Range("E16").Value = '  +1.47%  '
$ws.Range("E17").Value = '  +0.31%  '
$ws.Range("D18").Value = '79.76'
$ws.Range("E18").Value = '  +3.32%  '
$ws.Range("D19").Value = '0.000009034'
$ws.Range("E19").Value = '  +1.23%  '
$ws.Range("D20").Value = '1.002'
$ws.Range("E20").Value = '  +0.17%  '
$ws.Range("D21").Value = '15.37'
$ws.Range("E21").Value = '  +0.52%  '
$ws.Range("D22").Value = '27.698.52'
$ws.Range("E22").Value = '  +0.76%  '
$ws.Range("D23").Value = '4.971'
$ws.Range("E23").Value = '  +0.99%  '
$ws.Range("D24").Value = '10.39'
$ws.Range("E24").Value = '  -3.66%  '
$ws.Range("D25").Value = '2.115.64'
$ws.Range("E25").Value = '  +3.37%  '
$ws.Range("D26").Value = '1.963'
$ws.Range("E26").Value = '  -3.59%  '
$ws.Range("D27").Value = '153.89'
$ws.Range("E27").Value = '  +1.88%  '
$ws.Range("E28").Value = '  +3.50%  '
$ws.Range("D29").Value = '121.60'
$ws.Range("E29").Value = '  +9.65%  '
$ws.Range("D30").Value = '5.249'
$ws.Range("E30").Value = '  -1.26%  '
$ws.Range("D31").Value = '1.864'
$ws.Range("E31").Value = '  +9.06%  '
$ws.Range("D32").Value = '0.08922'
$ws.Range("E32").Value = '  +0.67%  '
$ws.Range("D33").Value = '0.7596'
$ws.Range("E33").Value = '  -1.87%  '
$ws.Range("D34").Value = '2.968'
$ws.Range("E34").Value = '  +3.85%  '
$ws.Range("E35").Value = '  +1.12%  '
$ws.Range("D36").Value = '1.100'
$ws.Range("E36").Value = '  +2.33%  '
$ws.Range("D37").Value = '0.05386'
$ws.Range("E37").Value = '  +0.07%  '
$ws.Range("D38").Value = '1.090'
$ws.Range("E38").Value = '  -0.15%  '
$ws.Range("E39").Value = '  +0.77%  '
$ws.Range("D40").Value = '2.811'
$ws.Range("E40").Value = '  -4.37%  '
$ws.Range("D41").Value = '0.5071'
$ws.Range("E41").Value = '  +0.49%  '
$ws.Range("D42").Value = '0.1652'
$ws.Range("E42").Value = '  +1.26%  '
$ws.Range("E43").Value = '  -0.18%  '
$ws.Range("D44").Value = '8.326'
$ws.Range("E44").Value = '  +1.75%  '
$ws.Range("D45").Value = '0.06541'
$ws.Range("E45").Value = '  -1.09%  '
$ws.Range("D46").Value = '10.35'
$ws.Range("E46").Value = '  +1.35%  '
$ws.Range("D47").Value = '105.00'
$ws.Range("E47").Value = '  -0.91%  '
$ws.Range("D48").Value = '0.4660'
$ws.Range("E48").Value = '  -0.76%  '
$ws.Range("E49").Value = '  +0.17%  '
$ws.Range("D50").Value = '1.617'
$ws.Range("E50").Value = '  -1.57%  '
$ws.Range("D51").Value = '1.758'
$ws.Range("E51").Value = '  -4.81%  '

$ws.Range("D5").Style = "Normal"
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").Style = "Normal"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").Style = "Normal"
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").Style = "Normal"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").Style = "Normal"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").Style = "Normal"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").Style = "Normal"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").Style = "Normal"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").Style = "Normal"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").Style = "Normal"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").Style = "Normal"
